# Weekly refresh of "Fruta / hortaliza" data for
# Hortaliza, Terminal La Palmera de La Serena - Ciboulette.
#
# The underlying source data re-shuffles which week's figures (Fecha plus
# the associated Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg) land on which worksheet row. Column D (Fecha)
# together with columns J, K, L, M and P move together as a unit from row
# to row; every other column (A, B, C, E-I, N, O, Q, R) stays untouched.
#
# Row => (D:Fecha, J:Volumen, K:Precio minimo, L:Precio maximo,
#         M:Precio promedio ponderado, P:Precio $/Kg)
$newValues = @{
    2  = @(45013, 1100, 2000, 2500, 2250, 750)
    3  = @(45028, 1000, 2000, 2500, 2250, 750)
    4  = @(44881, 500,  1900, 2000, 1950, 650)
    5  = @(44965, 1120, 2000, 2500, 2250, 750)
    6  = @(45006, 1100, 2000, 2500, 2250, 750)
    7  = @(45070, 800,  2000, 2500, 2250, 750)
    8  = @(44911, 700,  1800, 2000, 1900, 633)
    9  = @(45084, 900,  2000, 2500, 2250, 750)
    10 = @(44999, 1100, 2000, 2500, 2250, 750)
    11 = @(44827, 1200, 2000, 2500, 2250, 750)
    12 = @(44910, 1000, 1800, 2000, 1900, 633)
    13 = @(44985, 1000, 2000, 2500, 2250, 750)
    14 = @(44970, 800,  2000, 2500, 2250, 750)
    15 = @(45091, 800,  2000, 2500, 2250, 750)
    16 = @(45034, 1100, 2000, 2500, 2250, 750)
    17 = @(44883, 500,  1800, 2000, 1900, 633)
    18 = @(44685, 400,  1500, 2000, 1750, 583)
    19 = @(45020, 1200, 2000, 2500, 2250, 750)
    20 = @(45035, 1100, 2000, 2500, 2250, 750)
    21 = @(44953, 1000, 2000, 2500, 2250, 750)
    22 = @(44951, 800,  2000, 2500, 2250, 750)
    23 = @(44964, 1000, 2000, 2500, 2250, 750)
    24 = @(44978, 1000, 1800, 2000, 1900, 633)
    25 = @(45062, 1100, 2000, 2500, 2250, 750)
    26 = @(44992, 1040, 2000, 2500, 2250, 750)
    27 = @(44971, 1000, 2000, 2500, 2250, 750)
    28 = @(44848, 1000, 1500, 2000, 1750, 583)
    29 = @(45041, 1160, 2000, 2500, 2250, 750)
    30 = @(45077, 760,  2000, 2500, 2250, 750)
    31 = @(45007, 1160, 2000, 2500, 2250, 750)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P - Precio $/Kg
}
